$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point precision drift on A50's stored serial date value
$ws.Range("A50").Value = 44363.77004171528

# Append the newly retrieved row of data as row 51
$ws.Range("A51").Value = 44364.76966868654
$ws.Range("B51").Value = 78897
$ws.Range("C51").Value = 66382
$ws.Range("D51").Value = 3483
$ws.Range("E51").Value = 2118
$ws.Range("F51").Value = 1498
$ws.Range("G51").Value = 20778
$ws.Range("H51").Value = 1429
$ws.Range("I51").Value = 902
$ws.Range("J51").Value = 187

# Match the date-time number format used by the rest of column A
$ws.Range("A51").NumberFormat = $ws.Range("A50").NumberFormat
